$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.041240515573738
$ws.Cells.Item(2, 4).Value = 1.057568763589349
$ws.Cells.Item(2, 5).Value = 1.050415151866332
$ws.Cells.Item(2, 6).Value = 1.063400875303862
$ws.Cells.Item(2, 9).Value = 1.044027115974436
$ws.Cells.Item(2, 10).Value = 1.046322789755185
$ws.Cells.Item(2, 11).Value = 1.060303375646331
$ws.Cells.Item(2, 12).Value = 1.053169492995887
$ws.Cells.Item(2, 13).Value = 1.066119615497034
$ws.Cells.Item(2, 14).Value = 1.047808688131539

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.042126483843132
$ws.Cells.Item(3, 4).Value = 1.058236428656855
$ws.Cells.Item(3, 5).Value = 1.051188941203702
$ws.Cells.Item(3, 6).Value = 1.064223730400519
$ws.Cells.Item(3, 9).Value = 1.044226590187372
$ws.Cells.Item(3, 10).Value = 1.046855021992533
$ws.Cells.Item(3, 11).Value = 1.060785274120953
$ws.Cells.Item(3, 12).Value = 1.053755823764507
$ws.Cells.Item(3, 13).Value = 1.06675745439619
$ws.Cells.Item(3, 14).Value = 1.048341676199712

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.042700481853017
$ws.Cells.Item(4, 4).Value = 1.058668931598862
$ws.Cells.Item(4, 5).Value = 1.051690622553879
$ws.Cells.Item(4, 6).Value = 1.064757172358421
$ws.Cells.Item(4, 9).Value = 1.044354784745038
$ws.Cells.Item(4, 10).Value = 1.047199474458479
$ws.Cells.Item(4, 11).Value = 1.061096888984373
$ws.Cells.Item(4, 12).Value = 1.054135544309236
$ws.Cells.Item(4, 13).Value = 1.067170524239911
$ws.Cells.Item(4, 14).Value = 1.048686617827669

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.042941960850108
$ws.Cells.Item(5, 4).Value = 1.058850868653857
$ws.Cells.Item(5, 5).Value = 1.05190176412414
$ws.Cells.Item(5, 6).Value = 1.064981668758725
$ws.Cells.Item(5, 9).Value = 1.044408466677974
$ws.Cells.Item(5, 10).Value = 1.04734429603457
$ws.Cells.Item(5, 11).Value = 1.061227841230358
$ws.Cells.Item(5, 12).Value = 1.054295255314703
$ws.Cells.Item(5, 13).Value = 1.067344259906569
$ws.Cells.Item(5, 14).Value = 1.048831645067009

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.042982516175691
$ws.Cells.Item(6, 4).Value = 1.058881423247704
$ws.Cells.Item(6, 5).Value = 1.051937229435925
$ws.Cells.Item(6, 6).Value = 1.06501937655749
$ws.Cells.Item(6, 9).Value = 1.044417467713187
$ws.Cells.Item(6, 10).Value = 1.047368612985517
$ws.Cells.Item(6, 11).Value = 1.061249825678155
$ws.Cells.Item(6, 12).Value = 1.054322075940018
$ws.Cells.Item(6, 13).Value = 1.067373435599221
$ws.Cells.Item(6, 14).Value = 1.048855996550816

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.042703707840147
$ws.Cells.Item(7, 4).Value = 1.05867136220973
$ws.Cells.Item(7, 5).Value = 1.051693442917775
$ws.Cells.Item(7, 6).Value = 1.064760171159903
$ws.Cells.Item(7, 9).Value = 1.04435550287587
$ws.Cells.Item(7, 10).Value = 1.047201409518944
$ws.Cells.Item(7, 11).Value = 1.061098638975445
$ws.Cells.Item(7, 12).Value = 1.054137678076712
$ws.Cells.Item(7, 13).Value = 1.067172845387704
$ws.Cells.Item(7, 14).Value = 1.048688555636142

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.041539783281158
$ws.Cells.Item(8, 4).Value = 1.057794303830457
$ws.Cells.Item(8, 5).Value = 1.050676451975146
$ws.Cells.Item(8, 6).Value = 1.063678754886449
$ws.Cells.Item(8, 9).Value = 1.044094710703632
$ws.Cells.Item(8, 10).Value = 1.046502646425089
$ws.Cells.Item(8, 11).Value = 1.060466277233907
$ws.Cells.Item(8, 12).Value = 1.053367578022792
$ws.Cells.Item(8, 13).Value = 1.066335103660149
$ws.Cells.Item(8, 14).Value = 1.047988800218543

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.039494359077495
$ws.Cells.Item(9, 4).Value = 1.056252565455346
$ws.Cells.Item(9, 5).Value = 1.048892028742407
$ws.Cells.Item(9, 6).Value = 1.061780902847433
$ws.Cells.Item(9, 9).Value = 1.043628464538049
$ws.Cells.Item(9, 10).Value = 1.045271870662448
$ws.Cells.Item(9, 11).Value = 1.05935046193892
$ws.Cells.Item(9, 12).Value = 1.052013116743916
$ws.Cells.Item(9, 13).Value = 1.06486161162049
$ws.Cells.Item(9, 14).Value = 1.046756276613149

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.038134569439085
$ws.Cells.Item(10, 4).Value = 1.055227372248977
$ws.Cells.Item(10, 5).Value = 1.047707659118729
$ws.Cells.Item(10, 6).Value = 1.060520988123289
$ws.Cells.Item(10, 9).Value = 1.043313174156651
$ws.Cells.Item(10, 10).Value = 1.044451785968354
$ws.Cells.Item(10, 11).Value = 1.058605651542126
$ws.Cells.Item(10, 12).Value = 1.051111942582984
$ws.Cells.Item(10, 13).Value = 1.063881204552099
$ws.Cells.Item(10, 14).Value = 1.045935027304724

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.037546691660506
$ws.Cells.Item(11, 4).Value = 1.054784100050002
$ws.Cells.Item(11, 5).Value = 1.047196080621856
$ws.Cells.Item(11, 6).Value = 1.059976717246731
$ws.Cells.Item(11, 9).Value = 1.043175601260808
$ws.Cells.Item(11, 10).Value = 1.04409679826599
$ws.Cells.Item(11, 11).Value = 1.058282936687332
$ws.Cells.Item(11, 12).Value = 1.050722168718
$ws.Cells.Item(11, 13).Value = 1.063457153127976
$ws.Cells.Item(11, 14).Value = 1.045579535479104

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.03732846739318
$ws.Cells.Item(12, 4).Value = 1.054619547552739
$ws.Cells.Item(12, 5).Value = 1.047006248909298
$ws.Cells.Item(12, 6).Value = 1.059774745160128
$ws.Cells.Item(12, 9).Value = 1.043124343516253
$ws.Cells.Item(12, 10).Value = 1.043964958449153
$ws.Cells.Item(12, 11).Value = 1.058163036338691
$ws.Cells.Item(12, 12).Value = 1.050577456942695
$ws.Cells.Item(12, 13).Value = 1.063299713972572
$ws.Cells.Item(12, 14).Value = 1.045447508434606

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.037375270893108
$ws.Cells.Item(13, 4).Value = 1.054654840132675
$ws.Cells.Item(13, 5).Value = 1.047046959769335
$ws.Cells.Item(13, 6).Value = 1.059818060031946
$ws.Cells.Item(13, 9).Value = 1.043135345571569
$ws.Cells.Item(13, 10).Value = 1.043993237683911
$ws.Cells.Item(13, 11).Value = 1.058188756690851
$ws.Cells.Item(13, 12).Value = 1.050608495028313
$ws.Cells.Item(13, 13).Value = 1.063333481902541
$ws.Cells.Item(13, 14).Value = 1.04547582782912

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.037528650315779
$ws.Cells.Item(14, 4).Value = 1.054770496068216
$ws.Cells.Item(14, 5).Value = 1.04718038515994
$ws.Cells.Item(14, 6).Value = 1.059960018194974
$ws.Cells.Item(14, 9).Value = 1.043171367479181
$ws.Cells.Item(14, 10).Value = 1.044085899962647
$ws.Cells.Item(14, 11).Value = 1.058273026284402
$ws.Cells.Item(14, 12).Value = 1.05071020541046
$ws.Cells.Item(14, 13).Value = 1.063444137671495
$ws.Cells.Item(14, 14).Value = 1.045568621698921

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.037623170998159
$ws.Cells.Item(15, 4).Value = 1.054841768620663
$ws.Cells.Item(15, 5).Value = 1.047262618352564
$ws.Cells.Item(15, 6).Value = 1.060047509116625
$ws.Cells.Item(15, 9).Value = 1.043193540976752
$ws.Cells.Item(15, 10).Value = 1.044142994728892
$ws.Cells.Item(15, 11).Value = 1.058324943671539
$ws.Cells.Item(15, 12).Value = 1.050772881532765
$ws.Cells.Item(15, 13).Value = 1.063512325988124
$ws.Cells.Item(15, 14).Value = 1.045625797546282

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.038173604462873
$ws.Cells.Item(16, 4).Value = 1.055256804478908
$ws.Cells.Item(16, 5).Value = 1.047741637629043
$ws.Cells.Item(16, 6).Value = 1.060557136735789
$ws.Cells.Item(16, 9).Value = 1.043322282353842
$ws.Cells.Item(16, 10).Value = 1.044475347848494
$ws.Cells.Item(16, 11).Value = 1.058627064823566
$ws.Cells.Item(16, 12).Value = 1.051137819992064
$ws.Cells.Item(16, 13).Value = 1.063909357516353
$ws.Cells.Item(16, 14).Value = 1.045958622645436

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.03851912421462
$ws.Cells.Item(17, 4).Value = 1.055517319056885
$ws.Cells.Item(17, 5).Value = 1.048042452749323
$ws.Cells.Item(17, 6).Value = 1.060877156896246
$ws.Cells.Item(17, 9).Value = 1.043402757787668
$ws.Cells.Item(17, 10).Value = 1.044683855675909
$ws.Cells.Item(17, 11).Value = 1.058816523014474
$ws.Cells.Item(17, 12).Value = 1.05136685517648
$ws.Cells.Item(17, 13).Value = 1.064158532227359
$ws.Cells.Item(17, 14).Value = 1.046167426577883

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.038720748717707
$ws.Cells.Item(18, 4).Value = 1.055669334737402
$ws.Cells.Item(18, 5).Value = 1.048218034701532
$ws.Cells.Item(18, 6).Value = 1.061063942660521
$ws.Cells.Item(18, 9).Value = 1.043449596380978
$ws.Cells.Item(18, 10).Value = 1.044805485717127
$ws.Cells.Item(18, 11).Value = 1.05892701070029
$ws.Cells.Item(18, 12).Value = 1.051500489930991
$ws.Cells.Item(18, 13).Value = 1.064303916986415
$ws.Cells.Item(18, 14).Value = 1.046289229347709

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.038789512459157
$ws.Cells.Item(19, 4).Value = 1.055721178603491
$ws.Cells.Item(19, 5).Value = 1.048277924177024
$ws.Cells.Item(19, 6).Value = 1.061127652653287
$ws.Cells.Item(19, 9).Value = 1.043465549915712
$ws.Cells.Item(19, 10).Value = 1.044846960238212
$ws.Cells.Item(19, 11).Value = 1.058964680698139
$ws.Cells.Item(19, 12).Value = 1.05154606307871
$ws.Cells.Item(19, 13).Value = 1.064353497077847
$ws.Cells.Item(19, 14).Value = 1.046330762767371

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.038482044032202
$ws.Cells.Item(20, 4).Value = 1.055489361888267
$ws.Cells.Item(20, 5).Value = 1.048010165544929
$ws.Cells.Item(20, 6).Value = 1.060842808981665
$ws.Cells.Item(20, 9).Value = 1.043394134014791
$ws.Cells.Item(20, 10).Value = 1.044661483617461
$ws.Cells.Item(20, 11).Value = 1.058796198009438
$ws.Cells.Item(20, 12).Value = 1.051342277472264
$ws.Cells.Item(20, 13).Value = 1.064131793434973
$ws.Cells.Item(20, 14).Value = 1.046145022748545

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.037483480029967
$ws.Cells.Item(21, 4).Value = 1.0547364355426
$ws.Cells.Item(21, 5).Value = 1.047141089412042
$ws.Cells.Item(21, 6).Value = 1.059918209673347
$ws.Cells.Item(21, 9).Value = 1.043160764254671
$ws.Cells.Item(21, 10).Value = 1.044058612712638
$ws.Cells.Item(21, 11).Value = 1.058248211797431
$ws.Cells.Item(21, 12).Value = 1.050680252360736
$ws.Cells.Item(21, 13).Value = 1.063411550293315
$ws.Cells.Item(21, 14).Value = 1.045541295697886

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.036856452292072
$ws.Cells.Item(22, 4).Value = 1.054263611891123
$ws.Cells.Item(22, 5).Value = 1.046595774564875
$ws.Cells.Item(22, 6).Value = 1.059338002817485
$ws.Cells.Item(22, 9).Value = 1.043013127163073
$ws.Cells.Item(22, 10).Value = 1.043679671042236
$ws.Cells.Item(22, 11).Value = 1.05790350008564
$ws.Cells.Item(22, 12).Value = 1.050264402668764
$ws.Cells.Item(22, 13).Value = 1.062959124669056
$ws.Cells.Item(22, 14).Value = 1.045161815886848

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.037188774372715
$ws.Cells.Item(23, 4).Value = 1.054514209895385
$ws.Cells.Item(23, 5).Value = 1.046884750670041
$ws.Cells.Item(23, 6).Value = 1.059645474115204
$ws.Cells.Item(23, 9).Value = 1.043091478231382
$ws.Cells.Item(23, 10).Value = 1.043880544617215
$ws.Cells.Item(23, 11).Value = 1.058086254022992
$ws.Cells.Item(23, 12).Value = 1.050484814851065
$ws.Cells.Item(23, 13).Value = 1.063198923718437
$ws.Cells.Item(23, 14).Value = 1.045362974725344

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.038498798696706
$ws.Cells.Item(24, 4).Value = 1.055501994338241
$ws.Cells.Item(24, 5).Value = 1.04802475436737
$ws.Cells.Item(24, 6).Value = 1.060858328944681
$ws.Cells.Item(24, 9).Value = 1.043398031039838
$ws.Cells.Item(24, 10).Value = 1.044671592553984
$ws.Cells.Item(24, 11).Value = 1.05880538206653
$ws.Cells.Item(24, 12).Value = 1.05135338294801
$ws.Cells.Item(24, 13).Value = 1.064143875404323
$ws.Cells.Item(24, 14).Value = 1.046155146040918

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.040022482317837
$ws.Cells.Item(25, 4).Value = 1.056650685906393
$ws.Cells.Item(25, 5).Value = 1.049352428058283
$ws.Cells.Item(25, 6).Value = 1.062270613825486
$ws.Cells.Item(25, 9).Value = 1.043749789547207
$ws.Cells.Item(25, 10).Value = 1.045589984849382
$ws.Cells.Item(25, 11).Value = 1.059639097616246
$ws.Cells.Item(25, 12).Value = 1.05236296628172
$ws.Cells.Item(25, 13).Value = 1.065242213001379
$ws.Cells.Item(25, 14).Value = 1.047074842558716

